$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers
$ws.Range("A1").Value = "Cost Breakdown"
$ws.Range("B1").Value = "Cost_First Year"
$ws.Range("C1").Value = "Cost_Yearly"

# Update installation cost row
$ws.Range("A2").Value = "Convenient Installtion cost"
$ws.Range("B2").Value = 510

# Update the selection to B6
$ws.Range("B6").Select()
